$d = $word.ActiveDocument

$d.Content.Find.Execute("2024-06-19 Wednesday", $true, $false, $false, $false, $false, $true, 1, $false, "2024-06-20 Thursday", 2) | Out-Null
$d.Content.Find.Execute("50÷6=8, 2", $true, $false, $false, $false, $false, $true, 1, $false, "75÷9=8, 3", 2) | Out-Null
$d.Content.Find.Execute("16÷6=2, 4", $true, $false, $false, $false, $false, $true, 1, $false, "92÷2=46, 0", 2) | Out-Null
$d.Content.Find.Execute("41÷4=10, 1", $true, $false, $false, $false, $false, $true, 1, $false, "71÷2=35, 1", 2) | Out-Null
$d.Content.Find.Execute("64÷7=9, 1", $true, $false, $false, $false, $false, $true, 1, $false, "65÷7=9, 2", 2) | Out-Null
$d.Content.Find.Execute("84÷6=14, 0", $true, $false, $false, $false, $false, $true, 1, $false, "93÷9=10, 3", 2) | Out-Null
$d.Content.Find.Execute("83÷4=20, 3", $true, $false, $false, $false, $false, $true, 1, $false, "84÷6=14, 0", 2) | Out-Null
$d.Content.Find.Execute("17÷2=8, 1", $true, $false, $false, $false, $false, $true, 1, $false, "90÷9=10, 0", 2) | Out-Null
$d.Content.Find.Execute("33÷3=11, 0", $true, $false, $false, $false, $false, $true, 1, $false, "99÷9=11, 0", 2) | Out-Null
$d.Content.Find.Execute("46÷9=5, 1", $true, $false, $false, $false, $false, $true, 1, $false, "52÷2=26, 0", 2) | Out-Null
$d.Content.Find.Execute("89÷2=44, 1", $true, $false, $false, $false, $false, $true, 1, $false, "62÷3=20, 2", 2) | Out-Null
$d.Content.Find.Execute("66÷4=16, 2", $true, $false, $false, $false, $false, $true, 1, $false, "34÷8=4, 2", 2) | Out-Null
$d.Content.Find.Execute("24÷7=3, 3", $true, $false, $false, $false, $false, $true, 1, $false, "15÷8=1, 7", 2) | Out-Null
$d.Content.Find.Execute("67÷2=33, 1", $true, $false, $false, $false, $false, $true, 1, $false, "52÷9=5, 7", 2) | Out-Null
$d.Content.Find.Execute("31÷4=7, 3", $true, $false, $false, $false, $false, $true, 1, $false, "76÷4=19, 0", 2) | Out-Null
$d.Content.Find.Execute("48÷3=16, 0", $true, $false, $false, $false, $false, $true, 1, $false, "76÷3=25, 1", 2) | Out-Null
$d.Content.Find.Execute("78÷6=13, 0", $true, $false, $false, $false, $false, $true, 1, $false, "73÷4=18, 1", 2) | Out-Null
$d.Content.Find.Execute("26÷9=2, 8", $true, $false, $false, $false, $false, $true, 1, $false, "88÷8=11, 0", 2) | Out-Null
$d.Content.Find.Execute("74÷8=9, 2", $true, $false, $false, $false, $false, $true, 1, $false, "40÷3=13, 1", 2) | Out-Null
$d.Content.Find.Execute("47÷6=7, 5", $true, $false, $false, $false, $false, $true, 1, $false, "93÷7=13, 2", 2) | Out-Null
$d.Content.Find.Execute("46÷4=11, 2", $true, $false, $false, $false, $false, $true, 1, $false, "35÷2=17, 1", 2) | Out-Null
$d.Content.Find.Execute("46÷6=7, 4", $true, $false, $false, $false, $false, $true, 1, $false, "85÷2=42, 1", 2) | Out-Null
$d.Content.Find.Execute("15÷6=2, 3", $true, $false, $false, $false, $false, $true, 1, $false, "83÷7=11, 6", 2) | Out-Null
$d.Content.Find.Execute("80÷4=20, 0", $true, $false, $false, $false, $false, $true, 1, $false, "55÷6=9, 1", 2) | Out-Null
$d.Content.Find.Execute("27÷9=3, 0", $true, $false, $false, $false, $false, $true, 1, $false, "92÷3=30, 2", 2) | Out-Null
$d.Content.Find.Execute("99÷5=19, 4", $true, $false, $false, $false, $false, $true, 1, $false, "18÷8=2, 2", 2) | Out-Null
